# edit.ps1 - applies the HCAI module-description content edits described
# by the commit diff, using Word COM-interop (Find/Range) calls.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...iteratively testing his design and explicitly involving
#    stakeholders in the process." ->
#    "...iteratively testing his design and explicitly involving
#    stakeholders in the process to create a valuable holistic user
#    experience."
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("explicitly involving stakeholders in the process", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insRng = $d.Range($rng.End, $rng.End)
    $insRng.InsertAfter(" to create a valuable holistic user experience")
}

# ---------------------------------------------------------------------
# 2) "Specifically, the student can:" -> "Specifically, the students can:"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Specifically, the student can:", $true, $false, $false, $false, $false, $true, 1, $false, "Specifically, the students can:", 2)

# ---------------------------------------------------------------------
# 3) Item 1:
#    "1) develop a wireframe prototype using various iterative testing
#    methodologies, design techniques and processes" ->
#    "1) develop a wireframe prototype by applying various iterative
#    testing methodologies and, design- techniques and processes"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "1) develop a wireframe prototype using various iterative testing methodologies, design techniques and processes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1) develop a wireframe prototype by applying various iterative testing methodologies and, design- techniques and processes",
    2)

# ---------------------------------------------------------------------
# 4) Item 2:
#    "2) explicitly involving stakeholders in the design and testing of
#    application" ->
#    "2) design whilst explicitly involving stakeholders(users & clients)
#    in the design and testing of the prototype"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "2) explicitly involving stakeholders in the design and testing of application",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2) design whilst explicitly involving stakeholders(users & clients) in the design and testing of the prototype",
    2)

# ---------------------------------------------------------------------
# 5) Item 3:
#    "3) create a remote user test from his wireframe prototype" ->
#    "3) create a user test for their wireframe prototype"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "3) create a remote user test from his wireframe prototype",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3) create a user test for their wireframe prototype",
    2)

# ---------------------------------------------------------------------
# 6) Item 5:
#    "5) design for and communicate about various disruptive technology
#    risks" -> append " in the context of their application"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("5) design for and communicate about various disruptive technology risks", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insRng = $d.Range($rng.End, $rng.End)
    $insRng.InsertAfter(" in the context of their application")
}

Write-Output "Edits applied"
